$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new string-valued cells (A and F columns) in the exact
# order the original author entered them, so the shared-strings table
# gets rebuilt with the same index assignment as the target workbook. ---
$ws.Range("F15").Value = "(20,40)"
$ws.Range("F14").Value = "(10,20)"
$ws.Range("A15").Value = "EP80"
$ws.Range("A14").Value = "EP250"
$ws.Range("F16").Value = "(10,20,30)"
$ws.Range("F17").Value = "(10,20,30,40)"
$ws.Range("A16").Value = "EP170"
$ws.Range("A17").Value = "EP225"
$ws.Range("A18").Value = "EP160"
$ws.Range("F18").Value = "(10,20,30,40,50)"
$ws.Range("A19").Value = "EP400"
$ws.Range("F19").Value = "(20,40,60,80,100)"

# --- Row 14 (EP250 / (10,20)) ---
$ws.Range("B14").Value = 2293
$ws.Range("C14").Value = 0.001
$ws.Range("D14").Value = 0.5
$ws.Range("E14").Value = 250
$ws.Range("G14").Value = 0.099976
$ws.Range("H14").Value = 0.043611
$ws.Range("I14").Value = 0.106992
$ws.Range("J14").Value = 0.142608
$ws.Range("K14").Value = 0.182439

# --- Row 15 (EP80 / (20,40)) ---
$ws.Range("B15").Value = 2293
$ws.Range("C15").Value = 0.001
$ws.Range("D15").Value = 0.5
$ws.Range("E15").Value = 80
$ws.Range("G15").Value = 0.100019
$ws.Range("H15").Value = 0.043902
$ws.Range("I15").Value = 0.10612
$ws.Range("J15").Value = 0.141154
$ws.Range("K15").Value = 0.183021

# --- Row 16 (EP170 / (10,20,30)) ---
$ws.Range("B16").Value = 2293
$ws.Range("C16").Value = 0.001
$ws.Range("D16").Value = 0.5
$ws.Range("E16").Value = 170
$ws.Range("G16").Value = 0.099149
$ws.Range("H16").Value = 0.04521
$ws.Range("I16").Value = 0.110481
$ws.Range("J16").Value = 0.143189
$ws.Range("K16").Value = 0.189853

# --- Row 17 (EP225 / (10,20,30,40)) ---
$ws.Range("B17").Value = 2293
$ws.Range("C17").Value = 0.001
$ws.Range("D17").Value = 0.5
$ws.Range("E17").Value = 225
$ws.Range("G17").Value = 0.098771
$ws.Range("H17").Value = 0.04521
$ws.Range("I17").Value = 0.11019
$ws.Range("J17").Value = 0.143189
$ws.Range("K17").Value = 0.191888

# --- Row 18 (EP160 / (10,20,30,40,50)) ---
$ws.Range("B18").Value = 2293
$ws.Range("C18").Value = 0.001
$ws.Range("D18").Value = 0.5
$ws.Range("E18").Value = 160
$ws.Range("G18").Value = 0.098605
$ws.Range("H18").Value = 0.045792
$ws.Range("I18").Value = 0.110481
$ws.Range("J18").Value = 0.143916
$ws.Range("K18").Value = 0.190725

# --- Row 19 (EP400 / (20,40,60,80,100)) ---
$ws.Range("B19").Value = 2293
$ws.Range("C19").Value = 0.001
$ws.Range("D19").Value = 0.5
$ws.Range("E19").Value = 400
$ws.Range("G19").Value = 0.098605
$ws.Range("H19").Value = 0.044338
$ws.Range("I19").Value = 0.109173
$ws.Range("J19").Value = 0.144207
$ws.Range("K19").Value = 0.19058

# --- Column F width (widened to fit the longer "(...)" labels) ---
$ws.Columns.Item(6).ColumnWidth = 14.428571428571429

# --- Number format id renumbering cosmetic touch (best effort; keeps the
# same format code used by column G) ---
$ws.Range("G1:G19").NumberFormat = "0.000000_);[Red]\(0.000000\)"

# --- Final selection matches the author's last cursor position ---
$ws.Range("K18").Select()
